$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing the existing rows 4-8 down to 5-9
$ws.Rows(4).Insert()

# Populate the newly inserted row 4 with the new weekly record
$ws.Cells.Item(4,1).Value  = 11
$ws.Cells.Item(4,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(4,3).Value  = "Bíobío"
$ws.Cells.Item(4,4).Value  = 45133
$ws.Cells.Item(4,5).Value  = 8
$ws.Cells.Item(4,6).Value  = 100112035
$ws.Cells.Item(4,7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(4,8).Value  = "Sin especificar"
$ws.Cells.Item(4,9).Value  = "Primera"
$ws.Cells.Item(4,10).Value = 50
$ws.Cells.Item(4,11).Value = 22000
$ws.Cells.Item(4,12).Value = 22000
$ws.Cells.Item(4,13).Value = 22000
$ws.Cells.Item(4,14).Value = "$/malla 15 kilos"
$ws.Cells.Item(4,15).Value = "Provincia de Quillota"
$ws.Cells.Item(4,16).Value = 1467
$ws.Cells.Item(4,17).Value = 15
$ws.Cells.Item(4,18).Value = "Hortaliza"
